$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Two more Email/Password rows of login test data, following the same
# layout as the existing row (A2/B2).
$ws.Range("A3").Value = "sads@asdas.com"
$ws.Range("B3").Value = "sdw323"

$ws.Range("A4").Value = "sgy3@vcf.com"
$ws.Range("B4").Value = "yyguy234t"

# Mail-to hyperlinks on the new email cells, matching the existing A2
# hyperlink pattern (linking to the row's password value).
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:sdw323")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:yyguy234t")

# Row 2 (Email/Password) loses its highlight fill, keeping font/border.
$ws.Range("A2:B2").Interior.Pattern = -4142

# Leave the selection on the last cell that was filled in.
$ws.Range("B4").Select()
